$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data block for rows 2-10 (A:T), matching the updated TPM values
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$data = @(
    @("ECs","Cthrc1","Fzd3","ECs",1,0.3333333333333333,0.01838633333333333,0.055159,0.003339500866342531,0.003339500866342531,3,1,0.162165,0.486495,0.1006291402646046,0.1006291402646046,0.002981619745,0.026834577705,0.0003360511010929511,0.0003360511010929511),
    @("ECs","Cthrc1","Fzd3","FAPs",1,0.3333333333333333,0.01838633333333333,0.055159,0.003339500866342531,0.003339500866342531,3,1,0.307583,0.922749,0.1908661724170313,0.1908661724170313,0.005655323565666667,0.050897912091,0.0006373977481421588,0.0006373977481421588),
    @("ECs","Cthrc1","Fzd3","MuSCs",1,0.3333333333333333,0.01838633333333333,0.055159,0.003339500866342531,0.003339500866342531,3,1,1.141763333333333,3.42529,0.7085046873183641,0.7085046873183641,0.02099284123444445,0.18893557111,0.002366052017107421,0.002366052017107421),
    @("FAPs","Cthrc1","Fzd3","ECs",3,1,5.449420333333333,16.348261,0.9897755900704113,0.9897755900704112,3,1,0.162165,0.486495,0.1006291402646046,0.1006291402646046,0.8837052483549999,7.953347235194999,0.09960026668367718,0.09960026668367718),
    @("FAPs","Cthrc1","Fzd3","FAPs",3,1,5.449420333333333,16.348261,0.9897755900704113,0.9897755900704112,3,1,0.307583,0.922749,0.1908661724170313,0.1908661724170313,1.676149054387667,15.085341489489,0.188914678428548,0.188914678428548),
    @("FAPs","Cthrc1","Fzd3","MuSCs",3,1,5.449420333333333,16.348261,0.9897755900704113,0.9897755900704112,3,1,1.141763333333333,3.42529,0.7085046873183641,0.7085046873183641,6.221948324521112,55.99753492069001,0.7012606449581861,0.701260644958186),
    @("MuSCs","Cthrc1","Fzd3","ECs",1,0.3333333333333333,0.03790633333333333,0.113719,0.006884909063246366,0.006884909063246365,3,1,0.162165,0.486495,0.1006291402646046,0.1006291402646046,0.006147080544999999,0.055323724905,0.0006928224798344659,0.0006928224798344658),
    @("MuSCs","Cthrc1","Fzd3","FAPs",1,0.3333333333333333,0.03790633333333333,0.113719,0.006884909063246366,0.006884909063246365,3,1,0.307583,0.922749,0.1908661724170313,0.1908661724170313,0.01165934372566667,0.104934093531,0.001314096240341162,0.001314096240341162),
    @("MuSCs","Cthrc1","Fzd3","MuSCs",1,0.3333333333333333,0.03790633333333333,0.113719,0.006884909063246366,0.006884909063246365,3,1,1.141763333333333,3.42529,0.7085046873183641,0.7085046873183641,0.04328006150111112,0.38952055351,0.004877990343070738,0.004877990343070737)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
